# Feat : Load 함수 구현 - 드롭 프리펩 오류
# Fix the PrefabPath values in column H: strip the stray leading space and
# make each monster row (2-9 / 10) point at its own unique prefab path
# instead of several rows sharing "Monster/1/Red".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MonsterData")

$ws.Range("H2").Value  = "Monster/1/Blue"
$ws.Range("H3").Value  = "Monster/1/Green"
$ws.Range("H4").Value  = "Monster/1/Red"
$ws.Range("H5").Value  = "Monster/2/Red"
$ws.Range("H6").Value  = "Monster/3/Red"
$ws.Range("H7").Value  = "Monster/4/Red"
$ws.Range("H8").Value  = "Monster/5/Red"
$ws.Range("H9").Value  = "Monster/6/Red"
$ws.Range("H10").Value = "Monster/7/Red"

# Move the active selection to H6, matching where the author left off editing.
$ws.Range("H6").Select()
